# New crime data collected - weekly CompStat report refresh
# (Volume/Number label, reporting week dates, and the week/28-day/YTD/2-year
#  crime-count figures for rows 15-29 all roll forward to the new week.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header: volume/number label and "report covering the week" date range
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/15/2024  Through  1/21/2024"

# ---------------------------------------------------------------------
# Helper donors used to re-apply the "text" cell style (s=14) when a
# formerly-numeric cell becomes a blank-marker ("0" / "***.*") text cell.
# Column C already carries that style on rows 27-29.
# ---------------------------------------------------------------------

# Row 15 (MURDER - TRANSIT / Shooting Vic. block continues numbering)
$ws.Range("G15").Value = 2
$ws.Range("J15").Value = 2

# Row 16
$ws.Range("C16").NumberFormat = "#,##0"
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 300
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -10
$ws.Range("I16").Value = 5
$ws.Range("J16").Value = 8
$ws.Range("K16").Value = -37.5
$ws.Range("L16").Value = -44.444444444444
$ws.Range("M16").Value = -76.190476190476
$ws.Range("N16").Value = -90.74074074074

# Row 17
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -20
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 31
$ws.Range("H17").Value = -29.032258064516
$ws.Range("I17").Value = 13
$ws.Range("J17").Value = 23
$ws.Range("K17").Value = -43.478260869565
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 85.714285714285
$ws.Range("N17").Value = -7.142857142857

# Row 18
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 45.454545454545
$ws.Range("I18").Value = 13
$ws.Range("J18").Value = 9
$ws.Range("K18").Value = 44.444444444444
$ws.Range("L18").Value = 85.714285714285
$ws.Range("M18").Value = -18.75
$ws.Range("N18").Value = -84.33734939759

# Row 19
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -15.384615384615
$ws.Range("F19").Value = 35
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = -27.083333333333
$ws.Range("I19").Value = 21
$ws.Range("J19").Value = 39
$ws.Range("K19").Value = -46.153846153846
$ws.Range("L19").Value = -54.347826086956
$ws.Range("M19").Value = 31.25
$ws.Range("N19").Value = -41.666666666666

# Row 20
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = 42.857142857142
$ws.Range("F20").Value = 34
$ws.Range("G20").Value = 30
$ws.Range("H20").Value = 13.333333333333
$ws.Range("I20").Value = 29
$ws.Range("J20").Value = 22
$ws.Range("K20").Value = 31.818181818181
$ws.Range("L20").Value = 123.076923076923
$ws.Range("M20").Value = 31.818181818181
$ws.Range("N20").Value = -84.574468085106

# Row 21 (TOTAL row - bold)
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 117
$ws.Range("G21").Value = 132
$ws.Range("H21").Value = -11.363636363636
$ws.Range("I21").Value = 81
$ws.Range("J21").Value = 103
$ws.Range("K21").Value = -21.35922330097
$ws.Range("L21").Value = -8.988764044943
$ws.Range("M21").Value = -1.219512195121
$ws.Range("N21").Value = -78.795811518324

# Row 22 (Transit)
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("C22").Value = 1
$ws.Range("F22").NumberFormat = "#,##0"
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("I22").NumberFormat = "#,##0"
$ws.Range("I22").Value = 1
$ws.Range("L22").Value = -50
$ws.Range("M22").Value = 0

# Row 24 (Petit Larceny)
$ws.Range("C24").Value = 12
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = -36.842105263157
$ws.Range("F24").Value = 77
$ws.Range("G24").Value = 80
$ws.Range("H24").Value = -3.75
$ws.Range("I24").Value = 54
$ws.Range("J24").Value = 60
$ws.Range("K24").Value = -10
$ws.Range("L24").Value = -3.571428571428
$ws.Range("M24").Value = 58.823529411764

# Row 25 (Misd. Assault)
$ws.Range("C25").Value = 14
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = 27.272727272727
$ws.Range("F25").Value = 39
$ws.Range("G25").Value = 42
$ws.Range("H25").Value = -7.142857142857
$ws.Range("I25").Value = 32
$ws.Range("J25").Value = 29
$ws.Range("K25").Value = 10.344827586206
$ws.Range("L25").Value = 45.454545454545
$ws.Range("M25").Value = 39.130434782608

# Row 26 (UCR Rape*)
$ws.Range("G26").Value = 2
$ws.Range("J26").Value = 2

# Row 27 (Other Sex Crimes) - D/E become blank-marker text cells again
$ws.Range("C27").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("C27").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("C27").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = "***.*"

# Row 28 (Shooting Vic.)
$ws.Range("C28").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("C28").Copy()
$ws.Range("D28").PasteSpecial(-4122)

$ws.Range("C28").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = "***.*"

$ws.Range("L28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L28").Value = -100
$ws.Range("N28").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N28").Value = -100

# Row 29 (Shooting Inc.)
$ws.Range("C29").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("C29").Copy()
$ws.Range("D29").PasteSpecial(-4122)

$ws.Range("C29").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").Value = "***.*"

$ws.Range("L29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("L29").Value = -100
$ws.Range("N29").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("N29").Value = -100
